# Add the new reporting week "31 gennaio - 5 febbraio 2022" to all four
# sheets of the "Quadro Sintesi" workbook (one new data row per sheet).

$wb = $excel.ActiveWorkbook

$newLabel = "31 gennaio - 5 febbraio 2022"

# ---------------------------------------------------------------------
# Sheet "Classi" -> new row 5 (A1:M5)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Classi")

$ws1.Range("A5").Value = $newLabel
$ws1.Range("B5").Value = 6080
$ws1.Range("C5").Value = 8157
$ws1.Range("D5").Value = 0.745
$ws1.Range("D5").NumberFormat = "0.00%"
$ws1.Range("E5").Value = 375908
$ws1.Range("F5").Value = 279677
$ws1.Range("G5").Value = 0.744
$ws1.Range("G5").NumberFormat = "0.00%"
$ws1.Range("H5").Value = 246401
$ws1.Range("I5").Value = 32619
$ws1.Range("J5").Value = 0.881
$ws1.Range("J5").NumberFormat = "0.00%"
$ws1.Range("K5").Value = 0.117
$ws1.Range("K5").NumberFormat = "0.00%"
$ws1.Range("L5").Value = 33276
$ws1.Range("M5").Value = 0.119
$ws1.Range("M5").NumberFormat = "0.00%"

$ws1.Range("N5").Select()

# ---------------------------------------------------------------------
# Sheet "Alunni in presenza" -> new row 5 (A1:F5)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Alunni in presenza")

$ws2.Range("A5").Value = $newLabel
$ws2.Range("B5").Value = 7382391
$ws2.Range("B5").NumberFormat = "#,##0"
$ws2.Range("C5").Value = 5477689
$ws2.Range("C5").NumberFormat = "#,##0"
$ws2.Range("D5").Value = 0.742
$ws2.Range("D5").NumberFormat = "0.00%"
$ws2.Range("E5").Value = 4760985
$ws2.Range("E5").NumberFormat = "#,##0"
$ws2.Range("F5").Value = 0.869
$ws2.Range("F5").NumberFormat = "0.00%"

$ws2.Range("G12").Select()

# ---------------------------------------------------------------------
# Sheet "Alunni" -> new rows 14, 15, 16 (A1:F16)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Alunni")

$ws3.Range("A14").Value = $newLabel
$ws3.Range("B14").Value = "Infanzia"
$ws3.Range("C14").Value = 619532
$ws3.Range("C14").NumberFormat = "#,##0"
$ws3.Range("D14").Value = 519606
$ws3.Range("D14").NumberFormat = "#,##0"
$ws3.Range("E14").Value = 99926
$ws3.Range("E14").NumberFormat = "#,##0"
$ws3.Range("F14").Value = 0.161
$ws3.Range("F14").NumberFormat = "0.00%"

$ws3.Range("A15").Value = $newLabel
$ws3.Range("B15").Value = "Primaria"
$ws3.Range("C15").Value = 1707409
$ws3.Range("C15").NumberFormat = "#,##0"
$ws3.Range("D15").Value = 1442403
$ws3.Range("D15").NumberFormat = "#,##0"
$ws3.Range("E15").Value = 265006
$ws3.Range("E15").NumberFormat = "#,##0"
$ws3.Range("F15").Value = 0.155
$ws3.Range("F15").NumberFormat = "0.00%"

$ws3.Range("A16").Value = $newLabel
$ws3.Range("B16").Value = "Sec. 1° e 2° Grado"
$ws3.Range("C16").Value = 3150748
$ws3.Range("C16").NumberFormat = "#,##0"
$ws3.Range("D16").Value = 2798976
$ws3.Range("D16").NumberFormat = "#,##0"
$ws3.Range("E16").Value = 351772
$ws3.Range("E16").NumberFormat = "#,##0"
$ws3.Range("F16").Value = 0.112
$ws3.Range("F16").NumberFormat = "0.00%"

$ws3.Range("C14:E16").Select()

# ---------------------------------------------------------------------
# Sheet "Personale scolastico" -> new row 5 (A1:K6, row 6 pre-existing)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Personale scolastico")

$ws4.Range("A5").Value = $newLabel
$ws4.Range("B5").Value = 775867
$ws4.Range("B5").NumberFormat = "#,##0"
$ws4.Range("C5").Value = 572166
$ws4.Range("C5").NumberFormat = "#,##0"
$ws4.Range("D5").Value = 0.737
$ws4.Range("D5").NumberFormat = "0.00%"
$ws4.Range("E5").Value = 536000
$ws4.Range("E5").NumberFormat = "#,##0"
$ws4.Range("F5").Value = 0.937
$ws4.Range("F5").NumberFormat = "0.00%"
$ws4.Range("G5").Value = 204526
$ws4.Range("G5").NumberFormat = "#,##0"
$ws4.Range("H5").Value = 151847
$ws4.Range("H5").NumberFormat = "#,##0"
$ws4.Range("I5").Value = 0.742
$ws4.Range("I5").NumberFormat = "0.00%"
$ws4.Range("J5").Value = 144255
$ws4.Range("J5").NumberFormat = "#,##0"
$ws4.Range("K5").Value = 0.95
$ws4.Range("K5").NumberFormat = "0.00%"

$ws4.Activate()
$ws4.Range("J6").Select()
